$p = $ppt.ActivePresentation

# --- Slide 33: underline the first "expr" in "     => expr op expr op expr" ---
$s33 = $p.Slides.Item(33)
$shp33 = $s33.Shapes.Item(2)
$tr33 = $shp33.TextFrame.TextRange

# Locate "     => expr op expr op expr" (the paragraph that previously had no
# underlined "expr") and underline just the first "expr" that follows "=> ".
$fullText = $tr33.Text
$target = "     => expr op expr op expr"
$startIdx = $fullText.IndexOf($target)
while ($startIdx -ge 0) {
    # Character position (1-based) of the "expr" right after "     => "
    $exprStart = $startIdx + ("     => ").Length + 1
    $chars = $tr33.Characters($exprStart, 4)
    if ($chars.Text -eq "expr") {
        $chars.Font.Underline = -1
    }
    $startIdx = $fullText.IndexOf($target, $startIdx + 1)
}

# --- Slide 34: nudge "Group 70" up slightly ---
$s34 = $p.Slides.Item(34)
$grp = $s34.Shapes.Item(4)
if ($grp.Name -ne "Group 70") {
    for ($i = 1; $i -le $s34.Shapes.Count; $i++) {
        if ($s34.Shapes.Item($i).Name -eq "Group 70") {
            $grp = $s34.Shapes.Item($i)
            break
        }
    }
}
$grp.Top = 174.0
